# Update the LR-pairs TPM-derived values for Snca-Lag3 sheet.
# Columns: G, H, M, N, O, P, Q, R, S, T for rows 2-6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value2 = 0.017133
$ws.Range("H2").Value2 = 0.051399
$ws.Range("M2").Value2 = 7.085731
$ws.Range("N2").Value2 = 21.257193
$ws.Range("O2").Value2 = 0.1632021414334214
$ws.Range("P2").Value2 = 0.1632021414334214
$ws.Range("Q2").Value2 = 0.121399829223
$ws.Range("R2").Value2 = 1.092598463007
$ws.Range("S2").Value2 = 0.1632021414334214
$ws.Range("T2").Value2 = 0.1632021414334214

# Row 3 (M3/N3 unchanged)
$ws.Range("G3").Value2 = 0.017133
$ws.Range("H3").Value2 = 0.051399
$ws.Range("O3").Value2 = 0.2146673930709577
$ws.Range("P3").Value2 = 0.2146673930709577
$ws.Range("Q3").Value2 = 0.159682860958
$ws.Range("R3").Value2 = 1.437145748622
$ws.Range("S3").Value2 = 0.2146673930709577
$ws.Range("T3").Value2 = 0.2146673930709577

# Row 4
$ws.Range("G4").Value2 = 0.017133
$ws.Range("H4").Value2 = 0.051399
$ws.Range("M4").Value2 = 14.22029466666667
$ws.Range("N4").Value2 = 42.660884
$ws.Range("O4").Value2 = 0.3275290215525062
$ws.Range("P4").Value2 = 0.3275290215525062
$ws.Range("Q4").Value2 = 0.243636308524
$ws.Range("R4").Value2 = 2.192726776716
$ws.Range("S4").Value2 = 0.3275290215525062
$ws.Range("T4").Value2 = 0.3275290215525062

# Row 5
$ws.Range("G5").Value2 = 0.017133
$ws.Range("H5").Value2 = 0.051399
$ws.Range("M5").Value2 = 1.993361666666667
$ws.Range("N5").Value2 = 5.980085
$ws.Range("O5").Value2 = 0.04591211445245296
$ws.Range("P5").Value2 = 0.04591211445245296
$ws.Range("Q5").Value2 = 0.034152265435
$ws.Range("R5").Value2 = 0.307370388915
$ws.Range("S5").Value2 = 0.04591211445245296
$ws.Range("T5").Value2 = 0.04591211445245296

# Row 6
$ws.Range("G6").Value2 = 0.017133
$ws.Range("H6").Value2 = 0.051399
$ws.Range("M6").Value2 = 10.79731966666667
$ws.Range("N6").Value2 = 32.391959
$ws.Range("O6").Value2 = 0.2486893294906617
$ws.Range("P6").Value2 = 0.2486893294906617
$ws.Range("Q6").Value2 = 0.184990477849
$ws.Range("R6").Value2 = 1.664914300641
$ws.Range("S6").Value2 = 0.2486893294906617
$ws.Range("T6").Value2 = 0.2486893294906617
